# Add a "Personensteuer" column to the Steuerfuss table, inserted just
# before the existing last column ("JuristPerson" / BfsId is already the
# first column), shifting the former column J ("JuristPerson" values) to
# column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at J; this pushes the existing J column
# (JuristPerson header + its "0" data cells, style included) to K, and
# widens the sheet's used range from A1:J5 to A1:K5 automatically.
$ws.Range("J1").EntireColumn.Insert()

# Give the new column J its header text (added as a new shared string).
$ws.Range("J1").Value = "Personensteuer"

# Match the author's resulting selection in the saved file.
$ws.Range("J7").Select()
